$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    "B2" = 0.3770835877108993
    "C2" = 0.06204085455566144
    "D2" = 0.02428201727138912
    "E2" = 0.4169131172345146
    "F2" = 0.6745040499286858
    "K2" = 0.3628246890282298
    "N2" = 1.274553151604098
    "O2" = 2.286437771945032
    "B3" = 0.3350854292756082
    "C3" = 0.05623166161691984
    "D3" = 0.02260588466059232
    "E3" = 0.3637838837403251
    "F3" = 0.6703514289716281
    "K3" = 0.3177370021363402
    "N3" = 1.290348779401116
    "O3" = 2.286328194739156
    "B4" = 0.3093405411888739
    "C4" = 0.05264010522189722
    "D4" = 0.02156747538315074
    "E4" = 0.3312461650622112
    "F4" = 0.668253361976852
    "K4" = 0.290029613369768
    "N4" = 1.300534452545088
    "O4" = 2.287740803486457
    "B5" = 0.298860230531119
    "C5" = 0.05117033785475655
    "D5" = 0.02114201639175661
    "E5" = 0.3180062801570642
    "F5" = 0.6675118070990251
    "K5" = 0.2787330757165307
    "N5" = 1.304807596311643
    "O5" = 2.288687992941732
    "B6" = 0.2971206537965543
    "C6" = 0.05092591177815109
    "D6" = 0.02107123126406663
    "E6" = 0.3158089389431211
    "F6" = 0.6673955187059732
    "K6" = 0.2768569685403293
    "N6" = 1.305524538769181
    "O6" = 2.288867695942258
    "B7" = 0.3091991553053504
    "C7" = 0.05262030840498255
    "D7" = 0.02156174676385803
    "E7" = 0.3310675304697526
    "F7" = 0.6682429020565834
    "K7" = 0.2898772863599959
    "N7" = 1.300591586181331
    "O7" = 2.287752074055277
    "B8" = 0.3625940654347062
    "C8" = 0.06004298466108082
    "D8" = 0.02370602493751761
    "E8" = 0.398575785485292
    "F8" = 0.6729783912405196
    "K8" = 0.3472834846343744
    "N8" = 1.279898338295268
    "O8" = 2.286092475541579
    "B9" = 0.4676261129164629
    "C9" = 0.07440233283539044
    "D9" = 0.02783640828068457
    "E9" = 0.5317021362695442
    "F9" = 0.6858566540029472
    "K9" = 0.4596627964981792
    "N9" = 1.243186959970316
    "O9" = 2.294610304080408
    "B10" = 0.5449860507113726
    "C10" = 0.08483253962366177
    "D10" = 0.03082436763631335
    "E10" = 0.6300792176655534
    "F10" = 0.6975217766324562
    "K10" = 0.5421061580282469
    "N10" = 1.218576562687076
    "O10" = 2.308091560551418
    "B11" = 0.5802203928164431
    "C11" = 0.08955164447634445
    "D11" = 0.03217330660614692
    "E11" = 0.6749827859719204
    "F11" = 0.7033101054860111
    "K11" = 0.579585152635758
    "N11" = 1.207894173695782
    "O11" = 2.315803382655957
    "B12" = 0.5935686640664812
    "C12" = 0.09133494437793388
    "D12" = 0.03268260799116263
    "E12" = 0.69201026531357
    "F12" = 0.7055714866965985
    "K12" = 0.5937736860143445
    "N12" = 1.203922902070849
    "O12" = 2.318951489421892
    "B13" = 0.5906936235357705
    "C13" = 0.09095104494959116
    "D13" = 0.03257298859662683
    "E13" = 0.6883420216504419
    "F13" = 0.7050813653103916
    "K13" = 0.5907181129258561
    "N13" = 1.204774896588738
    "O13" = 2.318263344849385
    "B14" = 0.581318449251512
    "C14" = 0.08969843243639275
    "D14" = 0.03221523763146905
    "E14" = 0.6763831666734603
    "F14" = 0.7034947574093735
    "K14" = 0.5807525330474732
    "N14" = 1.207565972518395
    "O14" = 2.316057809504969
    "B15" = 0.5755766271189202
    "C15" = 0.08893068528364267
    "D15" = 0.03199590692435095
    "E15" = 0.6690611305115084
    "F15" = 0.7025319663315628
    "K15" = 0.574647802540909
    "N15" = 1.209285217770343
    "O15" = 2.314736544367577
    "B16" = 0.5426842346172407
    "C16" = 0.08452361595524849
    "D16" = 0.03073600136183785
    "E16" = 0.6271478550842744
    "F16" = 0.6971532065753934
    "K16" = 0.5396562850187365
    "N16" = 1.219285017447427
    "O16" = 2.307619413299363
    "B17" = 0.5225165380961698
    "C17" = 0.08181342782668821
    "D17" = 0.02996042984243275
    "E17" = 0.6014753923364538
    "F17" = 0.6939770441816648
    "K17" = 0.5181834819952655
    "N17" = 1.225551067461003
    "O17" = 2.303658255350484
    "B18" = 0.5109206748485633
    "C18" = 0.0802521871741817
    "D18" = 0.02951337413139044
    "E18" = 0.5867234126350667
    "F18" = 0.6921955384214868
    "K18" = 0.5058305713976381
    "N18" = 1.229203411001968
    "O18" = 2.301528493406948
    "B19" = 0.5069952241467774
    "C19" = 0.07972316468554652
    "D19" = 0.02936184354761906
    "E19" = 0.5817310227705832
    "F19" = 0.6916001325926686
    "K19" = 0.5016477004942033
    "N19" = 1.230448318999894
    "O19" = 2.30083289115251
    "B20" = 0.5246630061444364
    "C20" = 0.08210218187672069
    "D20" = 0.03004309111564396
    "E20" = 0.6042067938678031
    "F20" = 0.6943104583603485
    "K20" = 0.5204695415433491
    "N20" = 1.224879038154796
    "O20" = 2.304064543443786
    "B21" = 0.5840720110724078
    "C21" = 0.09006645640046429
    "D21" = 0.03232035906245301
    "E21" = 0.6798951215307767
    "F21" = 0.7039588958371326
    "K21" = 0.5836797755066243
    "N21" = 1.206744158572606
    "O21" = 2.316699440240058
    "B22" = 0.622932806178369
    "C22" = 0.09524985580911505
    "D22" = 0.03379985458066415
    "E22" = 0.7294994389327059
    "F22" = 0.7106696982456668
    "K22" = 0.6249681334511195
    "N22" = 1.195322883767864
    "O22" = 2.326285251598364
    "B23" = 0.6021891218030078
    "C23" = 0.09248537614871566
    "D23" = 0.03301103839235253
    "E23" = 0.70301152019114
    "F23" = 0.70705089919106
    "K23" = 0.6029340020307643
    "N23" = 1.201379154373726
    "O23" = 2.321047353218461
    "B24" = 0.5236925916059079
    "C24" = 0.08197164587939199
    "D24" = 0.03000572359317744
    "E24" = 0.6029719041563197
    "F24" = 0.6941595831670782
    "K24" = 0.5194360386356891
    "N24" = 1.225182707146136
    "O24" = 2.303880400991744
    "B25" = 0.4391779234185833
    "C25" = 0.0705387674002651
    "D25" = 0.02672713387346448
    "E25" = 0.4955964288447348
    "F25" = 0.6819867835137572
    "K25" = 0.4292823501930059
    "N25" = 1.252703924108974
    "O25" = 2.291040909878177
}

foreach ($cell in $values.Keys) {
    $ws.Range($cell).Value = $values[$cell]
}
